$d = $word.ActiveDocument

# --- Step 1: delete the ListBullet paragraph (old #21) "(Signature of the authorised party)" ---
# Find it by style + text to be robust to ordering, then delete its Range (removes the whole paragraph).
foreach ($para in $d.Paragraphs) {
    if ($para.Style.NameLocal -eq "List Bullet" -and $para.Range.Text.TrimEnd([char]13,[char]7) -eq '(Signature of the authorised party)') {
        $para.Range.Delete()
        break
    }
}

# --- Step 2: 1:1 text replacements (old paragraphs 1-20, 22-23 -> new paragraphs 1-22) ---
$d.Content.Find.Execute('OCR Results - Vollmacht  (2).pdf', $true, $false, $false, $false, $false, $true, 1, $false, 'OCR Results - diplome licence allemand.pdf', 2) | Out-Null
$d.Content.Find.Execute('Der Bevollmachtigte ist berechtigt, eine Untervolimacht, die den Umfang dieser Vollmacht nicht Uberschreiten', $true, $false, $false, $false, $false, $true, 1, $false, '——— |', 2) | Out-Null
$d.Content.Find.Execute('darf, zu erteilen und zu widerrufen. (The authorised party has the right to issue and to revoke a sub-', $true, $false, $false, $false, $false, $true, 1, $false, '. / Beeidigte Ubersetzung aus dem Arabischen', 2) | Out-Null
$d.Content.Find.Execute('authorisation, which must not exceed the scope of this authorisation.)', $true, $false, $false, $false, $false, $true, 1, $false, 'Republik Tunesien (Wappen der tunesischen Republik)', 2) | Out-Null
$d.Content.Find.Execute('Hinweis (Note)', $true, $false, $false, $false, $false, $true, 1, $false, 'Ministerium fiir Hochschulbildung und wissenschaftliche Forschung', 2) | Out-Null
$d.Content.Find.Execute('Auf die Méglichkelt der Verwendung des Musters fir die Untervolimacht als Aniage zu einer Volimacht nach § Abs. 1 AufenthG wird', $true, $false, $false, $false, $false, $true, 1, $false, 'Universitit von Manouba', 2) | Out-Null
$d.Content.Find.Execute('hingewiesen. (Please be advised that the form for the sub-authorisation may be used as an annex fo an authorisation in accordance with', $true, $false, $false, $false, $false, $true, 1, $false, 'Fakuiltit fiir Literaturwissenschaft, Kiinste und Humanwissenschaften', 2) | Out-Null
$d.Content.Find.Execute('Section (1) of the Residence Act)', $true, $false, $false, $false, $false, $true, 1, $false, 'Das Nationale Zeugnis der Fundamentalen Lizenz (Bachelor)', 2) | Out-Null
$d.Content.Find.Execute('Die Volimacht erlischt mit Abschiuss des beschieunigten Fachkrafteverfahrens. (The authorisation expires', $true, $false, $false, $false, $false, $true, 1, $false, '¥ Nach Vorlage des Erlasses Nr. 83 des Jahres 1986 vom . September 1986 aber die Errichtung der Fakultat der Literaturwissenschaft von Manouba,', 2) | Out-Null
$d.Content.Find.Execute('upon completion of the fast-track procedure for skilled workers.)', $true, $false, $false, $false, $false, $true, 1, $false, '¥ Nach Vorlage des Gesetzes Nr. 19 des Jahres 2008 vom 25. Februar 2008 aber das Hochschulwesen, insbesondere dessen Artikel 3,', 2) | Out-Null
$d.Content.Find.Execute('Da diese Volimacht meine rechtliche Méglichkeit selbst zu handeln nicht ausschlie&, bitte ich um direk-', $true, $false, $false, $false, $false, $true, 1, $false, '¥ Nach Vorlage des Erlasses Nr. 1932 des Jahres 2008 vom 02. November 1992 ober die Festiegung der Stelle, die die Unterzeichnung der wissenschaftichen nationalen Hochschulabschlasse zusttindig ist.', 2) | Out-Null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13,[char]7) -eq 'ten Kontakt zu mir, sofern dies zur Kiarung von Sachverhalten und zur Verfahrensbeschleunigung er-') {
        $para.Range.Text = '"Nach Vorlage des Erfasses Nr. 3123 des Jahres 2008 vom 22. September 2008 uber die Festlegung des allgemeinen Rakmens fbr das Studiensystem und der Bedingungen for den Erwerb des nationalen'
        break
    }
}
$d.Content.Find.Execute('forderlich erscheint. (Because this authorisation does not exclude my legal ability to take action myself, |', $true, $false, $false, $false, $false, $true, 1, $false, 'Hochschulabschlusses fur die 1.izenz in den verschiedenen Ausbildungsgebieten, Fachern, Studiengsngen und Fachrichtungen im LMD-System (Lizenz, Master, Doktor},', 2) | Out-Null
$d.Content.Find.Execute('ask that you contact me directly if this seems necessary in order to clarify circumstances and to speed up the', $true, $false, $false, $false, $false, $true, 1, $false, '¥ Und nach Vorlage der Beratungsprotokolle der Prifungskommissionen des Universitatsjahres 2015-2016,', 2) | Out-Null
$d.Content.Find.Execute('procedures.)', $true, $false, $false, $false, $false, $true, 1, $false, 'wird Frau/ Fraulein: Salma Njema (geboren am 01. 01. 1993 in Monastir, Nationalausweisnummer: 06935513)', 2) | Out-Null
$d.Content.Find.Execute('Vf (A ob. 9 a¢ / Fachkraft', $true, $false, $false, $false, $false, $true, 1, $false, 'Das Nationale Zeugnis der Fundamentalen Lizenz (Bachelor) in: Fachbereich: Sprachen und Literaturen Hauptfach: Deutsche Sprache, Literatur und Landeskunde mit dem Pradikat: (Ausreichend) erteilt.', 2) | Out-Null
$d.Content.Find.Execute('Ort, Datum (Place, date achtgeber', $true, $false, $false, $false, $false, $true, 1, $false, 'Manouba, den 02. 07. 2016', 2) | Out-Null
$d.Content.Find.Execute('(Signature of the granting authorisation)', $true, $false, $false, $false, $false, $true, 1, $false, 'Der Dekan: Habib Kozdoghli (Unterschrift: Unleserlich) - Dienstsiegel: (Fakultat flir Literaturwissenschaft, Kiinste und', 2) | Out-Null
$d.Content.Find.Execute('Bansin, 11.06.2025 Hotel', $true, $false, $false, $false, $false, $true, 1, $false, 'Humanwissenschaften-In der Mitte: Der Dekan).', 2) | Out-Null
$d.Content.Find.Execute('Ort, Datum (Place, date) Unterschrift Bevollmachtigte/Bevolimachtigter', $true, $false, $false, $false, $false, $true, 1, $false, 'Trockenes Dienstsiegel: (Ministerium fiir Hochschulbildung und wissenschaftliche Forschung- Universitat von Manouba —In der Mitte: Fakultat', 2) | Out-Null
$d.Content.Find.Execute('Beschleunigtes Fachkrafteverfahren — Volimacht nach § Abs. 1 AufenthG /', $true, $false, $false, $false, $false, $true, 1, $false, 'fiir Literaturwissenschaft, Kiinste und Humanwissenschaften von Manouba)', 2) | Out-Null
$d.Content.Find.Execute('Fast-track procedure for skilled workers — Authonsation in accordance with Section (1) of the German Residence Act Seite 3 von 3', $true, $false, $false, $false, $false, $true, 1, $false, 'Hinweis: Das vorliegende Diplom wird nur einmal ausgehindigt.', 2) | Out-Null

# --- Step 3: append new paragraphs (new #23-37) at the end of the document ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Auf der Riickseite:'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = '*Stempel des Ministeriums fiir Hochschulbildung und wissenschaftliche Forschung fir die Beglaubigung des Dokumentes: Beglaubigungsvermerk:'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Durchsicht erfolgte in der Generaldirektion fiir Hochschulbildung im Ministerium flir Hochschulbildung und wissenschaftliche Forschung. Hiermit bestitigen'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'wir die Echtheit der Unterschrift des Herm: Der Dekan ohne Verantwortung fiir den Inhalt des vorliegenden Dokuments, Beglaubigungsnummer: 3148, Ort'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'u. Datum: Tunis, den 22. 01. 2025, Beglaubigungsgebiihr: 5 Dinar, Vizedirektorin der privaten Hochschulbildung: Latifa Ben Abderrahmen Unterschrift (Unleserlich), Siegel des Ministeriums fir Hochschulbildung und wissenschaftliche Forschung (Republik Tunesien - Ministerium fiir'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Hochschulbildung und wissenschaftliche Forschung - In der Mitte: Wappen der tunesischen Republik)'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = '*Stempel des AuBenministeriums fiir die Beglaubigung des Dokumentes: Beglaubigungsvermerk: Durchsicht erfolgte im Ministerium fur auswartige'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Angelegenheiten. Hiermit bestatigen wir die Echtheit der Unterschrift der Frau: Latifa Ben Abderrahmen, i. A. des Ministers flr Hochschulbildung und'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'wissenschaftliche Forschung, Ort und Datum: Tunis, den 22. 01. 2025, Beglaubigungsgebihr: 5 Dinar, i. A. des Ministers fur auswirtige Angelegenheiten,'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'i, A. des Generaldirektors ftir konsularische Angelegenheiten: Hamida Labidi ~ Unterschrift (Unieserlich), Siegel des AuBenministeriums (Republik'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Tunesien - Ministerium ftir auswartige Angelegenheiten — In der Mitte: Wappen der tunesischen Republik)'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Der Ubersetzung ist eine Kopie des Dokuments angeheftet.'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'Die Richtigkeit und Vollstindigkeit vorstehender Ubersetzung des mir im Original vorgelegten :'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = 'und in arabischer Sprache abgefassten Dokuments wird hiermit bescheinigt. Tunis, den 26. 03. 2025 i'
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = '... |'
